$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell: new "time_taken" column, formatted like the other header cells (E1)
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("F1").Value = "time_taken"

# Data cells (plain, unstyled text values, matching source formatting)
$ws.Range("F2").Value = "2021-10-05 13:40:33.877619"
$ws.Range("F3").Value = "2021-10-05 13:40:33.877630"
$ws.Range("F4").Value = "2021-10-05 13:40:33.877633"
$ws.Range("F5").Value = "2021-10-05 13:40:33.877636"
$ws.Range("F6").Value = "2021-10-05 13:40:33.877639"
